# Update default (primary) header: remove the literal "Word" and "Hi"
# text runs that used to flank the two center/right paragraph tabs,
# leaving the ptab runs (and the watermark picture run) untouched.

$d = $word.ActiveDocument

# The document has a single section; its "default" header (rId8 /
# header2.xml) is the primary header, i.e. wdHeaderFooterPrimary (1).
$section = $d.Sections.Item(1)
$header = $section.Headers.Item(1)

# Remove "Word" (whole word match) by replacing the find range with
# an empty string - this deletes the entire run (incl. its rPr),
# matching the target XML exactly.
$rngWord = $header.Range
$rngWord.Find.Execute("Word", $true, $false, $false, $false, $false, $true, 1, $false, "", 2)

# Remove "Hi" the same way. Re-fetch the header range so the search
# starts from the top again.
$rngHi = $header.Range
$rngHi.Find.Execute("Hi", $true, $false, $false, $false, $false, $true, 1, $false, "", 2)
